# Auto-generated edit script: updates cached market-price / profit values
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, matching the latest
# scheduled-runner market-data refresh.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 15499.5
$ws.Range("J3").Value = 15499.5
$ws.Range("L3").Value = 15499.5
$ws.Range("N3").Value = -15727.5
# Row 9
$ws.Range("H9").Value = 84.666664
$ws.Range("I9").Value = 72.166664
$ws.Range("J9").Value = 109.666664
$ws.Range("K9").Value = 72.166664
$ws.Range("L9").Value = 109.666664
$ws.Range("M9").Value = 96.833336
$ws.Range("N9").Value = -447.666664
# Row 96
$ws.Range("H96").Value = 1579.3334
$ws.Range("J96").Value = 2349
$ws.Range("L96").Value = 7047
$ws.Range("N96").Value = -9793
# Row 102
$ws.Range("H102").Value = 15499.5
$ws.Range("J102").Value = 15499.5
$ws.Range("L102").Value = 15499.5
$ws.Range("N102").Value = -21989.5
# Row 116
$ws.Range("H116").Value = 3871.8572
$ws.Range("I116").Value = 3925.75
$ws.Range("K116").Value = 3925.75
$ws.Range("M116").Value = -483.75

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 937.3
$ws.Range("I2").Value = 830.3333
$ws.Range("J2").Value = 1900
$ws.Range("K2").Value = 830.3333
$ws.Range("L2").Value = 1900
$ws.Range("M2").Value = -717.3333
$ws.Range("N2").Value = -2126
# Row 61
$ws.Range("H61").Value = 4833.3335
$ws.Range("J61").Value = 5000
$ws.Range("L61").Value = 5000
$ws.Range("N61").Value = -5424
# Row 116
$ws.Range("H116").Value = 937.3
$ws.Range("I116").Value = 830.3333
$ws.Range("J116").Value = 1900
$ws.Range("K116").Value = 830.3333
$ws.Range("L116").Value = 1900
$ws.Range("M116").Value = 1463.6667
$ws.Range("N116").Value = -6488
# Row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
# Row 132
$ws.Range("H132").Value = 4203.385
$ws.Range("I132").Value = 4767.636
$ws.Range("J132").Value = 1100
$ws.Range("K132").Value = 14302.908
$ws.Range("L132").Value = 3300
$ws.Range("M132").Value = -11772.908
$ws.Range("N132").Value = -8360
# Row 136
$ws.Range("H136").Value = 4833.3335
$ws.Range("J136").Value = 5000
$ws.Range("L136").Value = 15000
$ws.Range("N136").Value = -20100

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 937.3
$ws.Range("I3").Value = 830.3333
$ws.Range("J3").Value = 1900
$ws.Range("K3").Value = 830.3333
$ws.Range("L3").Value = 1900
$ws.Range("M3").Value = -716.3333
$ws.Range("N3").Value = -2128
# Row 80
$ws.Range("H80").Value = 581.6
$ws.Range("I80").Value = 287
$ws.Range("J80").Value = 1170.8
$ws.Range("K80").Value = 287
$ws.Range("L80").Value = 1170.8
$ws.Range("M80").Value = 711
$ws.Range("N80").Value = -3166.8
# Row 83
$ws.Range("H83").Value = 581.6
$ws.Range("I83").Value = 287
$ws.Range("J83").Value = 1170.8
$ws.Range("K83").Value = 1435
$ws.Range("L83").Value = 5854
$ws.Range("M83").Value = 3557
$ws.Range("N83").Value = -15838
# Row 134
$ws.Range("H134").Value = 1224.5
$ws.Range("I134").Value = 1224.5
$ws.Range("K134").Value = 3673.5
$ws.Range("M134").Value = -1138.5
# Row 140
$ws.Range("H140").Value = 87593.336
$ws.Range("J140").Value = 87593.336
$ws.Range("L140").Value = 87593.336
$ws.Range("N140").Value = -97953.336

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 2298
$ws.Range("I99").Value = 2298
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2298
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -800
$ws.Range("N99").ClearContents()
# Row 126
$ws.Range("H126").Value = 2298
$ws.Range("I126").Value = 2298
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6894
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -4424
$ws.Range("N126").ClearContents()
# Row 132
$ws.Range("H132").Value = 1405
$ws.Range("I132").Value = 1405
$ws.Range("K132").Value = 4215
$ws.Range("M132").Value = -1685
# Row 138
$ws.Range("H138").Value = 77500
$ws.Range("J138").Value = 77500
$ws.Range("L138").Value = 77500
$ws.Range("N138").Value = -87780

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 39
$ws.Range("I2").Value = 25
$ws.Range("K2").Value = 150
$ws.Range("M2").Value = -37
# Row 56
$ws.Range("H56").Value = 3206.3333
$ws.Range("I56").Value = 3206.3333
$ws.Range("K56").Value = 3206.3333
$ws.Range("M56").Value = -2676.3333
# Row 58
$ws.Range("H58").Value = 2224.6667
$ws.Range("J58").Value = 2590
$ws.Range("L58").Value = 7770
$ws.Range("N58").Value = -8026
# Row 92
$ws.Range("H92").Value = 675.6667
$ws.Range("I92").Value = 720.8
$ws.Range("K92").Value = 2162.4
$ws.Range("M92").Value = -914.3999999999996
# Row 114
$ws.Range("H114").Value = 569.44446
$ws.Range("I114").Value = 258
$ws.Range("J114").Value = 1659.5
$ws.Range("K114").Value = 774
$ws.Range("L114").Value = 4978.5
$ws.Range("M114").Value = 2480
$ws.Range("N114").Value = -11486.5
# Row 129
$ws.Range("H129").Value = 1194
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 1194
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 3582
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -13582

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2028.4117
$ws.Range("I102").Value = 1626.2
$ws.Range("J102").Value = 5045
$ws.Range("K102").Value = 1626.2
$ws.Range("L102").Value = 5045
$ws.Range("M102").Value = -4.200000000000045
$ws.Range("N102").Value = -8289
# Row 113
$ws.Range("H113").Value = 2603.625
$ws.Range("I113").Value = 696.2727
$ws.Range("J113").Value = 6799.8
$ws.Range("K113").Value = 696.2727
$ws.Range("L113").Value = 6799.8
$ws.Range("M113").Value = 1473.7273
$ws.Range("N113").Value = -11139.8
# Row 122
$ws.Range("H122").Value = 5874.5
$ws.Range("I122").Value = 4599.2
$ws.Range("K122").Value = 13797.6
$ws.Range("M122").Value = -11347.6
# Row 132
$ws.Range("H132").Value = 1500
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 3000
$ws.Range("M132").Value = -470

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 5474.5
$ws.Range("I40").Value = 5587.1177
$ws.Range("J40").Value = 4836.3335
$ws.Range("K40").Value = 5587.1177
$ws.Range("L40").Value = 4836.3335
$ws.Range("M40").Value = -5451.1177
$ws.Range("N40").Value = -5108.3335
# Row 100
$ws.Range("H100").Value = 7984.7
$ws.Range("I100").Value = 1175
$ws.Range("J100").Value = 9687.125
$ws.Range("K100").Value = 1175
$ws.Range("L100").Value = 9687.125
$ws.Range("M100").Value = -634
$ws.Range("N100").Value = -10769.125
# Row 122
$ws.Range("H122").Value = 5000
$ws.Range("J122").Value = 5000
$ws.Range("L122").Value = 15000
$ws.Range("N122").Value = -19900

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 80
$ws.Range("H80").Value = 119999.5
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 119999.5
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 119999.5
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -121995.5
# Row 81
$ws.Range("H81").Value = 861.75
$ws.Range("I81").Value = 861.75
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1723.5
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -662.5
$ws.Range("N81").ClearContents()
# Row 83
$ws.Range("H83").Value = 119999.5
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 119999.5
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 359998.5
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -369982.5
# Row 84
$ws.Range("H84").Value = 861.75
$ws.Range("I84").Value = 861.75
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 8617.5
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -3313.5
$ws.Range("N84").ClearContents()
# Row 100
$ws.Range("H100").Value = 1320.5
$ws.Range("I100").Value = 884.6
$ws.Range("J100").Value = 3500
$ws.Range("K100").Value = 1769.2
$ws.Range("L100").Value = 7000
$ws.Range("M100").Value = -1228.2
$ws.Range("N100").Value = -8082
# Row 122
$ws.Range("H122").Value = 1171
$ws.Range("I122").Value = 668
$ws.Range("K122").Value = 2004
$ws.Range("M122").Value = 446
# Row 135
$ws.Range("H135").Value = 48607.5
$ws.Range("J135").Value = 48607.5
$ws.Range("L135").Value = 48607.5
$ws.Range("N135").Value = -58747.5
# Row 141
$ws.Range("H141").Value = 240355
$ws.Range("J141").Value = 153806.67
$ws.Range("L141").Value = 153806.67
$ws.Range("N141").Value = -164166.67

Write-Host "Applied scheduled market-data refresh to all leve profit sheets."
